$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New seed data (dynamic ticks / dynamic units rework) for
# "Non-residential" (col B) and "Residential" (col C)
$ws.Cells.Item(2, 2).Value = 3.5856
$ws.Cells.Item(2, 3).Value = 3.5358

$ws.Cells.Item(3, 2).Value = 3.5429
$ws.Cells.Item(3, 3).Value = 3.5928

$ws.Cells.Item(4, 2).Value = 3.6
$ws.Cells.Item(4, 3).Value = 3.55

$ws.Cells.Item(5, 2).Value = 3.5571
$ws.Cells.Item(5, 3).Value = 3.6072

$ws.Cells.Item(6, 2).Value = 3.6144
$ws.Cells.Item(6, 3).Value = 3.5642

$ws.Cells.Item(7, 2).Value = 3.5713
$ws.Cells.Item(7, 3).Value = 3.6216

$ws.Cells.Item(8, 2).Value = 3.5287
$ws.Cells.Item(8, 3).Value = 3.5784
